# Refresh the cryptocurrency price/volume snapshot (mirrors the GitHub
# Actions job that regenerates cryptos.xlsx on a schedule).
#
# Every value in columns B-E is stored as literal TEXT in the workbook
# (prices such as "37.111.87" use '.' as a thousands separator and are
# never real numbers; percentages keep their padding spaces). Excel's
# automatic "looks like a number" detection would otherwise coerce plain
# values like "253.19" or "1.00" into real numbers (losing the trailing
# zero / precision and picking up a different style), so every write
# temporarily forces a Text number format, assigns the value, then
# restores the cell's original style so no stray formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, in the order the rows appear in the sheet.
$updates = @(
    @{ Cell = "D2"; Value = "37.111.87" }
    @{ Cell = "E2"; Value = "  -0.08%  " }
    @{ Cell = "D3"; Value = "2.071.60" }
    @{ Cell = "E3"; Value = "  -0.98%  " }
    @{ Cell = "E4"; Value = "  -0.10%  " }
    @{ Cell = "D5"; Value = "253.19" }
    @{ Cell = "E5"; Value = "  +1.26%  " }
    @{ Cell = "D6"; Value = "0.676" }
    @{ Cell = "E6"; Value = "  +3.74%  " }
    @{ Cell = "D7"; Value = "61.68" }
    @{ Cell = "E7"; Value = "  +20.48%  " }
    @{ Cell = "D8"; Value = "0.999" }
    @{ Cell = "E8"; Value = "  -0.06%  " }
    @{ Cell = "B9"; Value = "Cardano" }
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada" }
    @{ Cell = "D9"; Value = "0.390" }
    @{ Cell = "E9"; Value = "  +6.01%  " }
    @{ Cell = "B10"; Value = "OKB" }
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb" }
    @{ Cell = "D10"; Value = "61.46" }
    @{ Cell = "E10"; Value = "  +1.57%  " }
    @{ Cell = "E11"; Value = "  +9.62%  " }
    @{ Cell = "E12"; Value = "  +3.00%  " }
    @{ Cell = "D13"; Value = "16.34" }
    @{ Cell = "E13"; Value = "  +7.71%  " }
    @{ Cell = "D14"; Value = "2.372.09" }
    @{ Cell = "E14"; Value = "  -0.74%  " }
    @{ Cell = "D15"; Value = "0.822" }
    @{ Cell = "E15"; Value = "  -0.24%  " }
    @{ Cell = "D16"; Value = "5.50" }
    @{ Cell = "E16"; Value = "  +9.00%  " }
    @{ Cell = "D17"; Value = "2.064.91" }
    @{ Cell = "E17"; Value = "  -1.14%  " }
    @{ Cell = "D18"; Value = "37.032.02" }
    @{ Cell = "E18"; Value = "  -0.12%  " }
    @{ Cell = "D19"; Value = "74.65" }
    @{ Cell = "E19"; Value = "  +3.72%  " }
    @{ Cell = "D20"; Value = "0.0₃0929" }
    @{ Cell = "E20"; Value = "  +13.17%  " }
    @{ Cell = "D21"; Value = "15.17" }
    @{ Cell = "E21"; Value = "  +13.71%  " }
    @{ Cell = "D22"; Value = "5.49" }
    @{ Cell = "E22"; Value = "  +5.69%  " }
    @{ Cell = "D23"; Value = "239.92" }
    @{ Cell = "E23"; Value = "  +0.20%  " }
    @{ Cell = "D24"; Value = "1.00" }
    @{ Cell = "E24"; Value = "  -0.09%  " }
    @{ Cell = "D25"; Value = "2.42" }
    @{ Cell = "E25"; Value = "  -0.02%  " }
    @{ Cell = "D26"; Value = "2.30" }
    @{ Cell = "E26"; Value = "  +15.49%  " }
    @{ Cell = "D27"; Value = "170.24" }
    @{ Cell = "E27"; Value = "  +0.69%  " }
    @{ Cell = "D28"; Value = "9.31" }
    @{ Cell = "E28"; Value = "  +2.17%  " }
    @{ Cell = "D29"; Value = "20.37" }
    @{ Cell = "E29"; Value = "  -0.54%  " }
    @{ Cell = "E30"; Value = "  +3.38%  " }
    @{ Cell = "D31"; Value = "4.79" }
    @{ Cell = "E31"; Value = "  +7.28%  " }
    @{ Cell = "E32"; Value = "  +6.78%  " }
    @{ Cell = "D33"; Value = "0.0639" }
    @{ Cell = "E33"; Value = "  +6.01%  " }
    @{ Cell = "D34"; Value = "4.41" }
    @{ Cell = "E34"; Value = "  +8.19%  " }
    @{ Cell = "D35"; Value = "0.0893" }
    @{ Cell = "E35"; Value = "  -2.29%  " }
    @{ Cell = "D36"; Value = "0.999" }
    @{ Cell = "E36"; Value = "  -0.11%  " }
    @{ Cell = "D37"; Value = "2.29" }
    @{ Cell = "E37"; Value = "  -1.08%  " }
    @{ Cell = "E38"; Value = "  -2.80%  " }
    @{ Cell = "D39"; Value = "0.111" }
    @{ Cell = "E39"; Value = "  +25.00%  " }
    @{ Cell = "D40"; Value = "1.37" }
    @{ Cell = "E40"; Value = "  +5.66%  " }
    @{ Cell = "D41"; Value = "18.10" }
    @{ Cell = "E41"; Value = "  +2.80%  " }
    @{ Cell = "D42"; Value = "0.0227" }
    @{ Cell = "E42"; Value = "  +1.97%  " }
    @{ Cell = "E43"; Value = "  +0.46%  " }
    @{ Cell = "D44"; Value = "98.66" }
    @{ Cell = "E44"; Value = "  +0.70%  " }
    @{ Cell = "D45"; Value = "4.43" }
    @{ Cell = "E45"; Value = "  +30.40%  " }
    @{ Cell = "E46"; Value = "  +2.96%  " }
    @{ Cell = "D47"; Value = "4.65" }
    @{ Cell = "E47"; Value = "  +15.79%  " }
    @{ Cell = "D48"; Value = "2.51" }
    @{ Cell = "E48"; Value = "  +10.79%  " }
    @{ Cell = "D49"; Value = "1.304.51" }
    @{ Cell = "E49"; Value = "  -0.53%  " }
    @{ Cell = "E50"; Value = "  -1.50%  " }
    @{ Cell = "D51"; Value = "6.91" }
    @{ Cell = "E51"; Value = "  +0.45%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $originalStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $u.Value
    $range.Style = $originalStyle
}
